# Swap the order of the "Recorded By" entries in column G from
# "System, dnasr281@gmail.com" to "dnasr281@gmail.com, System"
# for every row in the used range that matches the old value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldValue = "System, dnasr281@gmail.com"
$newValue = "dnasr281@gmail.com, System"

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # Column G = 7
    if ($cell.Value2 -eq $oldValue) {
        $cell.Value = $newValue
    }
}
